# [Kadastro App] Yeni kayit eklendi: 2974
#
# Appends a new record (row 41) to both the master "Kayitlar" sheet and the
# district-specific "Erdemli" sheet (the new record's Birim/district), then
# mirrors the same data onto the district sheet.

$wb = $excel.ActiveWorkbook

$recordNo  = "2974"
$tarih     = "2025-09-10"
$birim     = "Erdemli"
$parsel    = "1"
$is        = "ÇAP"
$personel  = "AYHAN KARADAYI (K.Teknisyeni)"

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 41

    # Columns A, B, D hold numeric-looking text ("2974", "2025-09-10", "1").
    # Force them to be stored as text (not auto-coerced to a number/date) by
    # pre-formatting the cell as Text, then restore the plain "Normal" style
    # once the value has been entered so no stray formatting is left behind.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $recordNo
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $tarih
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = $birim

    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $parsel
    $ws.Cells.Item($row, 4).Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $is

    $ws.Cells.Item($row, 6).Value = $personel
}
